$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the helper/test row (row 12, incrementing 1..39 formula row)
$ws.Rows.Item(12).Delete()

# Remove the beneficiary record for ERNESTINA / AAAE560802MZSLVR07 (row 11)
$ws.Rows.Item(11).Delete()

# Number the remaining beneficiary rows (2..10) sequentially in column AM
$ws.Range("AM2").Value = 1
$ws.Range("AM3").Value = 2
$ws.Range("AM4").Value = 3
$ws.Range("AM5").Value = 4
$ws.Range("AM6").Value = 5
$ws.Range("AM7").Value = 6
$ws.Range("AM8").Value = 7
$ws.Range("AM9").Value = 8
$ws.Range("AM10").Value = 9

# Fix up the hyperlink range that used to extend through the deleted row 11
$oldLink = $ws.Hyperlinks.Item(2)
$oldLink.Delete()
$ws.Hyperlinks.Add($ws.Range("AH3:AH10"), "mailto:example@hotmail.com", "", "", "example@hotmail.com")

# Update the view: scroll so column AI is at the left edge, and move the selection
$excel.ActiveWindow.ScrollColumn = 35
$ws.Range("AN9").Select()
